# New crime data collected: update the weekly CompStat report
# (volume/issue number, reporting week dates, and all weekly/28-day/
# year-to-date/2-year crime-complaint figures in the main table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -----------------------------------------------
# "Volume 32   Number  16" -> "...17"
$ws.Range("A8").Value = "Volume 32   Number  17"

# "Report Covering the Week  4/14/2025  Through  4/20/2025"
#   -> "...4/21/2025  Through  4/27/2025"
$ws.Range("C9").Value = "Report Covering the Week  4/21/2025  Through  4/27/2025"

# --- Row 14 (Murder) -----------------------------------------------------
$ws.Range("C14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 4
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 9
$ws.Range("J14").Value = 19
$ws.Range("K14").Value = -52.631578947368
$ws.Range("L14").Value = -57.142857142857
$ws.Range("M14").Value = -43.75
$ws.Range("N14").Value = -91.089108910891

# --- Row 15 (Rape) --------------------------------------------------------
$ws.Range("F15").Value = 13
$ws.Range("H15").Value = 18.181818181818
$ws.Range("I15").Value = 67
$ws.Range("J15").Value = 40
$ws.Range("K15").Value = 67.5
$ws.Range("L15").Value = 45.652173913043
$ws.Range("M15").Value = 13.559322033898
$ws.Range("N15").Value = -56.209150326797

# --- Row 16 (Robbery) ------------------------------------------------------
$ws.Range("C16").Value = 37
$ws.Range("E16").Value = -24.489795918367
$ws.Range("F16").Value = 112
$ws.Range("G16").Value = 170
$ws.Range("H16").Value = -34.117647058823
$ws.Range("I16").Value = 516
$ws.Range("J16").Value = 667
$ws.Range("K16").Value = -22.638680659670
$ws.Range("L16").Value = -6.352087114337
$ws.Range("M16").Value = -28.827586206896
$ws.Range("N16").Value = -82.334816843546

# --- Row 17 (Fel. Assault) --------------------------------------------------
$ws.Range("C17").Value = 56
$ws.Range("E17").Value = -13.846153846153
$ws.Range("F17").Value = 208
$ws.Range("G17").Value = 244
$ws.Range("H17").Value = -14.754098360655
$ws.Range("I17").Value = 843
$ws.Range("J17").Value = 956
$ws.Range("K17").Value = -11.820083682008
$ws.Range("L17").Value = -4.313280363223
$ws.Range("M17").Value = 52.166064981949
$ws.Range("N17").Value = -49.032648125755

# --- Row 18 (Burglary) ------------------------------------------------------
$ws.Range("C18").Value = 34
$ws.Range("D18").Value = 23
$ws.Range("E18").Value = 47.826086956521
$ws.Range("F18").Value = 107
$ws.Range("G18").Value = 96
$ws.Range("H18").Value = 11.458333333333
$ws.Range("I18").Value = 442
$ws.Range("J18").Value = 454
$ws.Range("K18").Value = -2.643171806167
$ws.Range("L18").Value = -17.537313432835
$ws.Range("M18").Value = -0.674157303370
$ws.Range("N18").Value = -87.251225843668

# --- Row 19 (Gr. Larceny) ----------------------------------------------------
$ws.Range("C19").Value = 132
$ws.Range("D19").Value = 137
$ws.Range("E19").Value = -3.649635036496
$ws.Range("F19").Value = 482
$ws.Range("G19").Value = 502
$ws.Range("H19").Value = -3.984063745019
$ws.Range("I19").Value = 1789
$ws.Range("J19").Value = 2003
$ws.Range("K19").Value = -10.683974038941
$ws.Range("L19").Value = -3.035230352303
$ws.Range("M19").Value = 28.243727598566
$ws.Range("N19").Value = -47.613469985358

# --- Row 20 (G.L.A.) ---------------------------------------------------------
$ws.Range("C20").Value = 20
$ws.Range("D20").Value = 17
$ws.Range("E20").Value = 17.647058823529
$ws.Range("F20").Value = 78
$ws.Range("G20").Value = 66
$ws.Range("H20").Value = 18.181818181818
$ws.Range("I20").Value = 251
$ws.Range("J20").Value = 288
$ws.Range("K20").Value = -12.847222222222
$ws.Range("L20").Value = -35.805626598465
$ws.Range("M20").Value = 60.897435897435
$ws.Range("N20").Value = -91.911053818885

# --- Row 21 (TOTAL) -----------------------------------------------------------
$ws.Range("C21").Value = 282
$ws.Range("D21").Value = 294
$ws.Range("E21").Value = -4.081632653061
$ws.Range("F21").Value = 1004
$ws.Range("G21").Value = 1093
$ws.Range("H21").Value = -8.142726440988
$ws.Range("I21").Value = 3917
$ws.Range("J21").Value = 4427
$ws.Range("K21").Value = -11.520216851140
$ws.Range("L21").Value = -8.288457035822
$ws.Range("M21").Value = 16.925373134328
$ws.Range("N21").Value = -73.558795733765

# --- Row 22 (Transit) ---------------------------------------------------------
$ws.Range("C22").Value = 6
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = 20
$ws.Range("G22").Value = 14
$ws.Range("H22").Value = -7.142857142857
$ws.Range("I22").Value = 71
$ws.Range("J22").Value = 88
$ws.Range("K22").Value = -19.318181818181
$ws.Range("L22").Value = -22.826086956521
$ws.Range("M22").Value = -4.054054054054

# --- Row 23 (Housing) ---------------------------------------------------------
$ws.Range("C23").Value = 20
$ws.Range("D23").Value = 25
$ws.Range("E23").Value = -20
$ws.Range("F23").Value = 97
$ws.Range("G23").Value = 96
$ws.Range("H23").Value = 1.041666666666
$ws.Range("I23").Value = 405
$ws.Range("J23").Value = 412
$ws.Range("K23").Value = -1.699029126213
$ws.Range("L23").Value = 5.46875
$ws.Range("M23").Value = 65.983606557377

# --- Row 24 (Petit Larceny) ----------------------------------------------------
$ws.Range("C24").Value = 252
$ws.Range("D24").Value = 269
$ws.Range("E24").Value = -6.319702602230
$ws.Range("F24").Value = 1027
$ws.Range("G24").Value = 1010
$ws.Range("H24").Value = 1.683168316831
$ws.Range("I24").Value = 4519
$ws.Range("J24").Value = 3973
$ws.Range("K24").Value = 13.742763654669
$ws.Range("L24").Value = 5.141926477431
$ws.Range("M24").Value = 69.632132132132

# --- Row 25 (Retail Theft) ------------------------------------------------------
$ws.Range("C25").Value = 145
$ws.Range("D25").Value = 138
$ws.Range("E25").Value = 5.072463768115
$ws.Range("F25").Value = 572
$ws.Range("G25").Value = 552
$ws.Range("H25").Value = 3.623188405797
$ws.Range("I25").Value = 2572
$ws.Range("J25").Value = 2128
$ws.Range("K25").Value = 20.864661654135
$ws.Range("L25").Value = 7.615062761506

# --- Row 26 (Misd. Assault) -----------------------------------------------------
$ws.Range("C26").Value = 133
$ws.Range("D26").Value = 91
$ws.Range("E26").Value = 46.153846153846
$ws.Range("F26").Value = 412
$ws.Range("G26").Value = 373
$ws.Range("H26").Value = 10.455764075067
$ws.Range("I26").Value = 1474
$ws.Range("J26").Value = 1485
$ws.Range("K26").Value = -0.740740740740
$ws.Range("L26").Value = 5.738880918220
$ws.Range("M26").Value = -4.471808165910

# --- Row 27 (UCR Rape*) -----------------------------------------------------------
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = -40
$ws.Range("F27").Value = 18
$ws.Range("G27").Value = 18
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 82
$ws.Range("J27").Value = 66
$ws.Range("K27").Value = 24.242424242424
$ws.Range("L27").Value = 0

# --- Row 28 (Other Sex Crimes) -----------------------------------------------------
$ws.Range("C28").Value = 14
$ws.Range("D28").Value = 9
$ws.Range("E28").Value = 55.555555555555
$ws.Range("F28").Value = 47
$ws.Range("G28").Value = 39
$ws.Range("H28").Value = 20.512820512820
$ws.Range("I28").Value = 175
$ws.Range("J28").Value = 174
$ws.Range("K28").Value = 0.574712643678
$ws.Range("L28").Value = -5.405405405405

# --- Row 29 (Shooting Vic.) ---------------------------------------------------------
$ws.Range("C29").Value = 5
$ws.Range("E29").Value = 150
$ws.Range("F29").Value = 8
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 25
$ws.Range("J29").Value = 31
$ws.Range("K29").Value = -19.354838709677
$ws.Range("L29").Value = -47.916666666666
$ws.Range("M29").Value = -52.830188679245
$ws.Range("N29").Value = -89.583333333333

# --- Row 30 (Shooting Inc.) ----------------------------------------------------------
$ws.Range("C30").Value = 4
$ws.Range("E30").Value = 100
$ws.Range("F30").Value = 7
$ws.Range("G30").Value = 7
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 22
$ws.Range("J30").Value = 25
$ws.Range("K30").Value = -12
$ws.Range("L30").Value = -51.111111111111
$ws.Range("M30").Value = -55.102040816326
$ws.Range("N30").Value = -90.134529147982

# --- Row 31 (Hate Crimes) --------------------------------------------------------------
$ws.Range("D31").Value = 3
$ws.Range("F31").Value = 5
$ws.Range("H31").Value = -61.538461538461
$ws.Range("I31").Value = 24
$ws.Range("J31").Value = 39
$ws.Range("K31").Value = -38.461538461538
$ws.Range("L31").Value = 4.347826086956

# --- Row 33 (Traffic Fatalities) --------------------------------------------------------
# C33 switches from the text placeholder "0" to a real number (2), so the
# number format has to move from General to the same #,##0 format used by
# the other numeric cells in this row (style used by F33/G33/I33/J33).
$ws.Range("C33").NumberFormat = "#,##0"
$ws.Range("C33").Value = 2
$ws.Range("F33").Value = 3
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 50
$ws.Range("I33").Value = 7
$ws.Range("K33").Value = 40
$ws.Range("L33").Value = 0
